$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Helper to force a literal text value into a cell even when the text
# looks like a percentage (e.g. "41.9%"), which Excel would otherwise
# auto-convert into a numeric percent value. We build the text through a
# formula (so it is never type-inferred from typed input) and then use
# Copy / PasteSpecial(values) to bake it back down to a static value
# while preserving the cell's existing style.
function Set-LiteralText($rangeAddress, [string]$text) {
    $cell = $ws.Range($rangeAddress)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# Row 3 - recorded-by list reordered
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 4 - recorded-by list reordered and attendance count updated
$ws.Range("G4").Value = "rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H4").Value = "157/221"

# Row 10 - average attendance % statistic
Set-LiteralText "L10" "41.9%"

# Row 12 - recorded-by list reordered
$ws.Range("G12").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"

# Row 15 - avg attendance % statistic
Set-LiteralText "S15" "49.6%"

# Row 16 - avg attendance % statistic
Set-LiteralText "S16" "36.2%"

# Row 25 - recorded-by list reordered
$ws.Range("G25").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 26 - recorded-by list reordered and attendance count updated
$ws.Range("G26").Value = "rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("H26").Value = "122/246"

# Row 34 - recorded-by list reordered
$ws.Range("G34").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"

# Row 41 - recorded-by list reordered
$ws.Range("G41").Value = "Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
